{"js": "// Update the answer key table: each of the 25 filled-in cells (5 columns x\n// 5 \"problem\" rows, interleaved with 3 blank rows each) gets its\n// \"NN\u00f7N=Q, R\" text replaced with a new equation string. Cells are addressed\n// by their (row, col) position inside the single table in the document\n// body, which is robust regardless of run/paragraph boundaries.\nconst changes = [\n  { row: 0, col: 0, oldText: \"48\u00f74=12, 0\", newText: \"39\u00f79=4, 3\" },\n  { row: 0, col: 1, oldText: \"53\u00f77=7, 4\", newText: \"73\u00f72=36, 1\" },\n  { row: 0, col: 2, oldText: \"29\u00f77=4, 1\", newText: \"44\u00f74=11, 0\" },\n  { row: 0, col: 3, oldText: \"54\u00f72=27, 0\", newText: \"35\u00f77=5, 0\" },\n  { row: 0, col: 4, oldText: \"34\u00f73=11, 1\", newText: \"64\u00f77=9, 1\" },\n  { row: 4, col: 0, oldText: \"80\u00f72=40, 0\", newText: \"36\u00f73=12, 0\" },\n  { row: 4, col: 1, oldText: \"92\u00f78=11, 4\", newText: \"57\u00f77=8, 1\" },\n  { row: 4, col: 2, oldText: \"23\u00f73=7, 2\", newText: \"90\u00f74=22, 2\" },\n  { row: 4, col: 3, oldText: \"37\u00f79=4, 1\", newText: \"21\u00f76=3, 3\" },\n  { row: 4, col: 4, oldText: \"78\u00f77=11, 1\", newText: \"57\u00f76=9, 3\" },\n  { row: 8, col: 0, oldText: \"75\u00f76=12, 3\", newText: \"91\u00f73=30, 1\" },\n  { row: 8, col: 1, oldText: \"58\u00f76=9, 4\", newText: \"54\u00f75=10, 4\" },\n  { row: 8, col: 2, oldText: \"29\u00f79=3, 2\", newText: \"31\u00f72=15, 1\" },\n  { row: 8, col: 3, oldText: \"92\u00f75=18, 2\", newText: \"83\u00f76=13, 5\" },\n  { row: 8, col: 4, oldText: \"95\u00f77=13, 4\", newText: \"59\u00f78=7, 3\" },\n  { row: 12, col: 0, oldText: \"92\u00f78=11, 4\", newText: \"72\u00f73=24, 0\" },\n  { row: 12, col: 1, oldText: \"71\u00f74=17, 3\", newText: \"21\u00f74=5, 1\" },\n  { row: 12, col: 2, oldText: \"45\u00f75=9, 0\", newText: \"23\u00f78=2, 7\" },\n  { row: 12, col: 3, oldText: \"96\u00f72=48, 0\", newText: \"92\u00f72=46, 0\" },\n  { row: 12, col: 4, oldText: \"70\u00f75=14, 0\", newText: \"84\u00f79=9, 3\" },\n  { row: 16, col: 0, oldText: \"74\u00f73=24, 2\", newText: \"91\u00f79=10, 1\" },\n  { row: 16, col: 1, oldText: \"78\u00f75=15, 3\", newText: \"41\u00f72=20, 1\" },\n  { row: 16, col: 2, oldText: \"48\u00f79=5, 3\", newText: \"19\u00f74=4, 3\" },\n  { row: 16, col: 3, oldText: \"88\u00f73=29, 1\", newText: \"78\u00f74=19, 2\" },\n  { row: 16, col: 4, oldText: \"21\u00f78=2, 5\", newText: \"47\u00f74=11, 3\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load the current text of every target cell first (one round trip) so we\n// can sanity-check against the expected \"old\" value before writing.\nconst cells = changes.map((c) => table.getCell(c.row, c.col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < changes.length; i++) {\n  const cell = cells[i];\n  const { newText } = changes[i];\n  // Skip cells that already hold the target text so the script is\n  // idempotent/safe if run more than once (e.g. against its own output).\n  if (cell.value === newText) {\n    continue;\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer key table: each of the 25 filled-in cells (5 columns x\n# 5 \"problem\" rows, interleaved with 3 blank rows each) gets its\n# \"NN\u00f7N=Q, R\" text replaced with a new equation string. Cells are addressed\n# by their 1-based (Row, Col) position inside the single table in the\n# document, which is robust regardless of run/paragraph boundaries.\n$changes = @(\n    @{ Row = 1; Col = 1; OldText = \"48\u00f74=12, 0\"; NewText = \"39\u00f79=4, 3\" },\n    @{ Row = 1; Col = 2; OldText = \"53\u00f77=7, 4\"; NewText = \"73\u00f72=36, 1\" },\n    @{ Row = 1; Col = 3; OldText = \"29\u00f77=4, 1\"; NewText = \"44\u00f74=11, 0\" },\n    @{ Row = 1; Col = 4; OldText = \"54\u00f72=27, 0\"; NewText = \"35\u00f77=5, 0\" },\n    @{ Row = 1; Col = 5; OldText = \"34\u00f73=11, 1\"; NewText = \"64\u00f77=9, 1\" },\n    @{ Row = 5; Col = 1; OldText = \"80\u00f72=40, 0\"; NewText = \"36\u00f73=12, 0\" },\n    @{ Row = 5; Col = 2; OldText = \"92\u00f78=11, 4\"; NewText = \"57\u00f77=8, 1\" },\n    @{ Row = 5; Col = 3; OldText = \"23\u00f73=7, 2\"; NewText = \"90\u00f74=22, 2\" },\n    @{ Row = 5; Col = 4; OldText = \"37\u00f79=4, 1\"; NewText = \"21\u00f76=3, 3\" },\n    @{ Row = 5; Col = 5; OldText = \"78\u00f77=11, 1\"; NewText = \"57\u00f76=9, 3\" },\n    @{ Row = 9; Col = 1; OldText = \"75\u00f76=12, 3\"; NewText = \"91\u00f73=30, 1\" },\n    @{ Row = 9; Col = 2; OldText = \"58\u00f76=9, 4\"; NewText = \"54\u00f75=10, 4\" },\n    @{ Row = 9; Col = 3; OldText = \"29\u00f79=3, 2\"; NewText = \"31\u00f72=15, 1\" },\n    @{ Row = 9; Col = 4; OldText = \"92\u00f75=18, 2\"; NewText = \"83\u00f76=13, 5\" },\n    @{ Row = 9; Col = 5; OldText = \"95\u00f77=13, 4\"; NewText = \"59\u00f78=7, 3\" },\n    @{ Row = 13; Col = 1; OldText = \"92\u00f78=11, 4\"; NewText = \"72\u00f73=24, 0\" },\n    @{ Row = 13; Col = 2; OldText = \"71\u00f74=17, 3\"; NewText = \"21\u00f74=5, 1\" },\n    @{ Row = 13; Col = 3; OldText = \"45\u00f75=9, 0\"; NewText = \"23\u00f78=2, 7\" },\n    @{ Row = 13; Col = 4; OldText = \"96\u00f72=48, 0\"; NewText = \"92\u00f72=46, 0\" },\n    @{ Row = 13; Col = 5; OldText = \"70\u00f75=14, 0\"; NewText = \"84\u00f79=9, 3\" },\n    @{ Row = 17; Col = 1; OldText = \"74\u00f73=24, 2\"; NewText = \"91\u00f79=10, 1\" },\n    @{ Row = 17; Col = 2; OldText = \"78\u00f75=15, 3\"; NewText = \"41\u00f72=20, 1\" },\n    @{ Row = 17; Col = 3; OldText = \"48\u00f79=5, 3\"; NewText = \"19\u00f74=4, 3\" },\n    @{ Row = 17; Col = 4; OldText = \"88\u00f73=29, 1\"; NewText = \"78\u00f74=19, 2\" },\n    @{ Row = 17; Col = 5; OldText = \"21\u00f78=2, 5\"; NewText = \"47\u00f74=11, 3\" }\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nforeach ($change in $changes) {\n    $cell = $table.Cell($change.Row, $change.Col)\n    # Cell.Range.Text always carries a trailing cell-mark (\"`r`a\"); compare\n    # against the stored value with that suffix stripped so the check works\n    # whether or not the marker is included.\n    $current = $cell.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($current -eq $change.NewText) {\n        continue\n    }\n    $cell.Range.Text = $change.NewText\n}\n"}
